# Auto-generated edit script applying numeric "want-to-go" count updates
# plus a content correction for two rows on the "展览" (Exhibition) sheet.
$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 289
$ws.Range("F5").Value = 2943
$ws.Range("F7").Value = 238
$ws.Range("F10").Value = 6907
$ws.Range("F12").Value = 73
$ws.Range("F13").Value = 350
$ws.Range("F16").Value = 1113
$ws.Range("F17").Value = 2233
$ws.Range("F18").Value = 1480
$ws.Range("F19").Value = 651
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 1110
$ws.Range("F22").Value = 125
$ws.Range("F23").Value = 178
$ws.Range("F24").Value = 344
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 1715
$ws.Range("F27").Value = 1693
$ws.Range("F28").Value = 1031
$ws.Range("F31").Value = 1221
$ws.Range("F32").Value = 140
$ws.Range("F33").Value = 583
$ws.Range("F34").Value = 32
$ws.Range("F36").Value = 426
$ws.Range("F37").Value = 13
$ws.Range("F38").Value = 2466
$ws.Range("F39").Value = 2720
$ws.Range("F43").Value = 15
$ws.Range("F44").Value = 24
$ws.Range("F48").Value = 154
$ws.Range("F49").Value = 415

# Row 41/42 content correction: a new con (夜蓝诗·恋与深空同人only) was
# inserted ahead of the existing "无限流同人only" entry, and the old
# "第七届燃梦BACG PRO...(取消)" row-42 entry was replaced.
$ws.Range("C41").Value = "上海·夜蓝诗·恋与深空同人only"
$ws.Range("D41").Value = "莫干山路50号 M50创意园"
$ws.Range("E41").Value = "2024.09.15 11:00-09.15 21:00"
$ws.Range("F41").Value = 2
$ws.Range("G41").Value = 88
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=90729"
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202408/dBiBf2Ac1723543844923.jpeg"
$ws.Range("C42").Value = "上海·无限流同人only"
$ws.Range("D42").Value = "呼青路158号 交运智慧湾科创园25号楼"
$ws.Range("E42").Value = "2024.09.15 10:00-09.15 17:00"
$ws.Range("F42").Value = 185
$ws.Range("G42").Value = 89
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=90108"
$ws.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202407/P3XVrcMn1722407440627.jpeg"

# ---- 演出 (Performance) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 182
$ws.Range("F8").Value = 212
$ws.Range("F12").Value = 182
$ws.Range("F19").Value = 47
$ws.Range("F20").Value = 50
$ws.Range("F23").Value = 470

# ---- 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1680
$ws.Range("F8").Value = 2723
$ws.Range("F9").Value = 1008
$ws.Range("F10").Value = 923
$ws.Range("F12").Value = 264
$ws.Range("F13").Value = 1455
$ws.Range("F14").Value = 7339

# ---- 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2944
$ws.Range("F6").Value = 1680
$ws.Range("F8").Value = 2723
$ws.Range("F9").Value = 6907
$ws.Range("F10").Value = 1008
$ws.Range("F12").Value = 350
$ws.Range("F14").Value = 264
$ws.Range("F15").Value = 1113
$ws.Range("F16").Value = 2233
$ws.Range("F17").Value = 1480
$ws.Range("F18").Value = 115
$ws.Range("F19").Value = 182
$ws.Range("F20").Value = 1110
$ws.Range("F22").Value = 1715
$ws.Range("F26").Value = 1221
$ws.Range("F27").Value = 140
$ws.Range("F29").Value = 583
$ws.Range("F30").Value = 32
$ws.Range("F31").Value = 50
$ws.Range("F34").Value = 470
$ws.Range("F35").Value = 426
$ws.Range("F37").Value = 13
$ws.Range("F38").Value = 2466
$ws.Range("F39").Value = 2720
$ws.Range("F42").Value = 15
$ws.Range("F43").Value = 24
$ws.Range("F48").Value = 415

